$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "W:\Data\Forecast\Tools\forecast_git\create_forecast_basic\current"
$ws.Range("B3").Select()
